$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33

$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 44610
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100112030
$ws.Cells.Item($row, 7).Value = "Poroto granado"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 80
$ws.Cells.Item($row, 11).Value = 30000
$ws.Cells.Item($row, 12).Value = 30000
$ws.Cells.Item($row, 13).Value = 30000
$ws.Cells.Item($row, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 1200
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"

$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(32, 4).NumberFormat

